$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "designs" (sheet1): fix rows 12-14 (drop the duplicate/extra
# "Output price" row, shift values up, delete the now-redundant row 14)
# -----------------------------------------------------------------
$wsDesigns = $wb.Worksheets.Item("designs")

$wsDesigns.Range("C12").Value = "Output efficiency"

$wsDesigns.Range("C13").Value = "Output price"
$wsDesigns.Range("E13").Value = 0
$wsDesigns.Range("F13").Value = '$/mi'

$wsDesigns.Rows("14").Delete()

# -----------------------------------------------------------------
# Sheet "parameters" (sheet2): add a new "Tractor lifetime" parameter row
# -----------------------------------------------------------------
$wsParameters = $wb.Worksheets.Item("parameters")

$wsParameters.Range("A34").Value = "Class 8 Diesel Tractor"
$wsParameters.Range("B34").Value = "Reference"
$wsParameters.Range("C34").Value = "Tractor lifetime"
$wsParameters.Range("D34").Value = 32
$wsParameters.Range("E34").Value = 10
$wsParameters.Range("F34").Value = "year"
$wsParameters.Range("G34").Value = "Vehicle lifetime before being retired"

# -----------------------------------------------------------------
# Sheet "results" (sheet3): insert a new "Output VMT" row at the top
# of the metrics table
# -----------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("results")

$wsResults.Rows("2").Insert()
$wsResults.Range("A2").Value = "Class 8 Diesel Tractor"
$wsResults.Range("B2").Value = "Output"
$wsResults.Range("C2").Value = "VMT"
$wsResults.Range("D2").Value = "mile/yr"
$wsResults.Range("E2").Value = "Yearly vehicle miles travelled (scale)"

# -----------------------------------------------------------------
# View-state / selection fixups to mirror the saved workbook state
# -----------------------------------------------------------------

# designs: selection now on row 12, no longer the active tab
$wsDesigns.Activate()
$wsDesigns.Rows("12").Select()

# parameters: selection moved to D14
$wsParameters.Activate()
$wsParameters.Range("D14").Select()

# indices (sheet4): selection moved to B11:B15
$wsIndices = $wb.Worksheets.Item("indices")
$wsIndices.Activate()
$wsIndices.Range("B11:B15").Select()

# results: becomes the active tab, selection on E3
$wsResults.Activate()
$wsResults.Range("E3").Select()

$wb.Save()
